$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-19, columns A-F (collected date, published date, org, title, ko title, link)
$data = @{
  2 = @("2026-01-30", "2026-01-29", "OpenAI", "Inside OpenAI’s in-house data agent", "OpenAI의 사내 데이터 에이전트 내부", "https://openai.com/index/inside-our-in-house-data-agent")
  3 = @("2026-01-30", "2026-01-29", "OpenAI", "Taisei Corporation shapes the next generation of talent with ChatGPT", "Taisei Corporation은 ChatGPT를 통해 차세대 인재를 양성합니다.", "https://openai.com/index/taisei")
  4 = @("2026-01-30", "2026-01-29", "OpenAI", "Retiring GPT-4o, GPT-4.1, GPT-4.1 mini, and OpenAI o4-mini in ChatGPT", "ChatGPT에서 GPT-4o, GPT-4.1, GPT-4.1 mini 및 OpenAI o4-mini 사용 중지", "https://openai.com/index/retiring-gpt-4o-and-older-models")
  5 = @("2026-01-29", "2026-01-28", "OpenAI", "The next chapter for AI in the EU", "EU AI의 다음 장", "https://openai.com/index/the-next-chapter-for-ai-in-the-eu")
  6 = @("2026-01-29", "2026-01-28", "OpenAI", "EMEA Youth & Wellbeing Grant", "EMEA 청소년 및 웰빙 보조금", "https://openai.com/index/emea-youth-and-wellbeing-grant")
  7 = @("2026-01-29", "2026-01-28", "OpenAI", "Keeping your data safe when an AI agent clicks a link", "AI 에이전트가 링크를 클릭할 때 데이터를 안전하게 유지", "https://openai.com/index/ai-agent-link-safety")
  8 = @("2026-01-28", "2026-01-27", "OpenAI", "PVH reimagines the future of fashion with OpenAI", "PVH는 OpenAI를 통해 패션의 미래를 재구상합니다.", "https://openai.com/index/pvh-future-of-fashion")
  9 = @("2026-01-28", "2026-01-27", "OpenAI", "Powering tax donations with AI powered personalized recommendations", "AI 기반 맞춤형 추천으로 세금 기부 지원", "https://openai.com/index/trustbank")
  10 = @("2026-01-28", "2026-01-27", "OpenAI", "Introducing Prism", "프리즘 소개", "https://openai.com/index/introducing-prism")
  11 = @("2026-01-28", "2026-01-26", "OpenAI", "How Indeed uses AI to help evolve the job search", "인디드가 AI를 활용하여 구직 활동을 발전시키는 방법", "https://openai.com/index/indeed-maggie-hulce")
  12 = @("2026-01-28", "2026-01-23", "OpenAI", "Unrolling the Codex agent loop", "Codex 에이전트 루프 풀기", "https://openai.com/index/unrolling-the-codex-agent-loop")
  13 = @("2026-01-28", "2026-01-22", "OpenAI", "Scaling PostgreSQL to power 800 million ChatGPT users", "8억 명의 ChatGPT 사용자를 지원하기 위해 PostgreSQL 확장", "https://openai.com/index/scaling-postgresql")
  14 = @("2026-01-28", "2026-01-22", "OpenAI", "Inside Praktika's conversational approach to language learning", "Praktika의 언어 학습에 대한 대화식 접근 방식 살펴보기", "https://openai.com/index/praktika")
  15 = @("2026-01-28", "2026-01-22", "OpenAI", "Inside GPT-5 for Work: How Businesses Use GPT-5", "업무용 GPT-5 내부: 기업이 GPT-5를 사용하는 방법", "https://openai.com/business/guides-and-resources/chatgpt-usage-and-adoption-patterns-at-work")
  16 = @("2026-01-29", "2026-01-21", "OpenAI", "How countries can end the capability overhang", "국가가 역량 과잉을 어떻게 끝낼 수 있는가", "https://openai.com/index/how-countries-can-end-the-capability-overhang")
  17 = @("2026-01-28", "2026-01-21", "OpenAI", "How Higgsfield turns simple ideas into cinematic social videos", "Higgsfield가 단순한 아이디어를 영화 같은 소셜 비디오로 바꾸는 방법", "https://openai.com/index/higgsfield")
  18 = @("2026-01-28", "2026-01-21", "OpenAI", "Introducing Edu for Countries", "국가별 교육 소개", "https://openai.com/index/edu-for-countries")
  19 = @("2026-01-29", "2026-01-20", "OpenAI", "Horizon 1000: Advancing AI for primary healthcare", "Horizon 1000: 1차 의료를 위한 AI 발전", "https://openai.com/index/horizon-1000")
}

# Columns A and B hold date-like text ("YYYY-MM-DD"); force text format so Excel
# does not auto-convert them into date serial numbers, matching the source inlineStr text.
$dateRange = $ws.Range("A2:B19")
$dateRange.NumberFormat = "@"

foreach ($r in $data.Keys) {
  $row = $data[$r]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item([int]$r, $c + 1).Value = $row[$c]
  }
}

# Clear the temporary text-format styling so cells keep the default (unstyled) look
# the source file uses for data rows.
$dateRange.ClearFormats()